# Scheduled runner update: refresh market-derived profit figures (columns
# H-N) across the per-job Leve profit tables (ALC, ARM, BSM, CRP, CUL,
# GSM, LTW, WVR). Values are written directly (no formulas are used in
# this workbook), mirroring the upstream data-refresh job.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 2119.2307
$ws.Range("I19").Value = 2566.2
$ws.Range("K19").Value = 2566.2
$ws.Range("M19").Value = -2391.2

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 6349.8667
$ws.Range("I51").Value = 4399.8
$ws.Range("J51").Value = 6739.88
$ws.Range("K51").Value = 4399.8
$ws.Range("L51").Value = 6739.88
$ws.Range("M51").Value = -3915.8
$ws.Range("N51").Value = -7707.88

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H74").Value = 5605.385
$ws.Range("I74").Value = 5605.385
$ws.Range("K74").Value = 5605.385
$ws.Range("M74").Value = -4669.385

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H77").Value = 5605.385
$ws.Range("I77").Value = 5605.385
$ws.Range("K77").Value = 28026.925
$ws.Range("M77").Value = -23346.925

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 6478.095
$ws.Range("I116").Value = 5514.2
$ws.Range("J116").Value = 8887.833000000001
$ws.Range("K116").Value = 5514.2
$ws.Range("L116").Value = 8887.833000000001
$ws.Range("M116").Value = -2072.2
$ws.Range("N116").Value = -15771.833

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 17262.6
$ws.Range("I137").Value = 8772.25
$ws.Range("K137").Value = 26316.75
$ws.Range("M137").Value = -23766.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2530.1396
$ws.Range("I2").Value = 1907.069
$ws.Range("J2").Value = 3820.7856
$ws.Range("K2").Value = 1907.069
$ws.Range("L2").Value = 3820.7856
$ws.Range("M2").Value = -1794.069
$ws.Range("N2").Value = -4046.7856

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H3").Value = 0
$ws.Range("I3").Value = 0
$ws.Range("K3").Value = 0
$ws.Range("M3").Value = $null

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6498.515
$ws.Range("I32").Value = 5853.484
$ws.Range("K32").Value = 5853.484
$ws.Range("M32").Value = -5566.484

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 14222.462
$ws.Range("J61").Value = 23199.6
$ws.Range("L61").Value = 23199.6
$ws.Range("N61").Value = -23623.6

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 30688.223
$ws.Range("I74").Value = 33274.25
$ws.Range("J74").Value = 10000
$ws.Range("K74").Value = 33274.25
$ws.Range("L74").Value = 10000
$ws.Range("M74").Value = -32400.25
$ws.Range("N74").Value = -11748

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 30688.223
$ws.Range("I77").Value = 33274.25
$ws.Range("J77").Value = 10000
$ws.Range("K77").Value = 166371.25
$ws.Range("L77").Value = 50000
$ws.Range("M77").Value = -162003.25
$ws.Range("N77").Value = -58736

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 1105.037
$ws.Range("I102").Value = 1105.037
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 1105.037
$ws.Range("L102").Value = 0
$ws.Range("M102").Value = 516.963
$ws.Range("N102").Value = $null

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H116").Value = 2530.1396
$ws.Range("I116").Value = 1907.069
$ws.Range("J116").Value = 3820.7856
$ws.Range("K116").Value = 1907.069
$ws.Range("L116").Value = 3820.7856
$ws.Range("M116").Value = 386.931
$ws.Range("N116").Value = -8408.785599999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 4175.3213
$ws.Range("I132").Value = 2545.6
$ws.Range("J132").Value = 8249.625
$ws.Range("K132").Value = 7636.799999999999
$ws.Range("L132").Value = 24748.875
$ws.Range("M132").Value = -5106.799999999999
$ws.Range("N132").Value = -29808.875

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 14222.462
$ws.Range("J136").Value = 23199.6
$ws.Range("L136").Value = 69598.79999999999
$ws.Range("N136").Value = -74698.79999999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2530.1396
$ws.Range("I3").Value = 1907.069
$ws.Range("J3").Value = 3820.7856
$ws.Range("K3").Value = 1907.069
$ws.Range("L3").Value = 3820.7856
$ws.Range("M3").Value = -1793.069
$ws.Range("N3").Value = -4048.7856

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1458.8695
$ws.Range("I94").Value = 927.8
$ws.Range("J94").Value = 4999.3335
$ws.Range("K94").Value = 927.8
$ws.Range("L94").Value = 4999.3335
$ws.Range("M94").Value = -476.8
$ws.Range("N94").Value = -5901.3335

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 8246.273999999999
$ws.Range("I134").Value = 3995.7646
$ws.Range("J134").Value = 16747.295
$ws.Range("K134").Value = 11987.2938
$ws.Range("L134").Value = 50241.88499999999
$ws.Range("M134").Value = -9452.293799999999
$ws.Range("N134").Value = -55311.88499999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H3").Value = 17000
$ws.Range("I3").Value = 14000
$ws.Range("J3").Value = 20000
$ws.Range("K3").Value = 14000
$ws.Range("L3").Value = 20000
$ws.Range("N3").Value = -20226
$ws.Range("M3").Value = -13887

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 462.7
$ws.Range("I22").Value = 223.89473
$ws.Range("K22").Value = 223.89473
$ws.Range("M22").Value = 126.10527

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 6420.857
$ws.Range("I58").Value = 3754.6
$ws.Range("J58").Value = 8844.727999999999
$ws.Range("K58").Value = 3754.6
$ws.Range("L58").Value = 8844.727999999999
$ws.Range("M58").Value = -3551.6
$ws.Range("N58").Value = -9250.727999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 3500
$ws.Range("I105").Value = 3250
$ws.Range("J105").Value = 4000
$ws.Range("K105").Value = 3250
$ws.Range("L105").Value = 4000
$ws.Range("M105").Value = -1503
$ws.Range("N105").Value = -7494

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 34726.65
$ws.Range("I132").Value = 23212.285
$ws.Range("K132").Value = 69636.855
$ws.Range("M132").Value = -67106.855

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 5285.4595
$ws.Range("I134").Value = 4024.963
$ws.Range("J134").Value = 8688.799999999999
$ws.Range("K134").Value = 12074.889
$ws.Range("L134").Value = 26066.4
$ws.Range("M134").Value = -9539.889000000001
$ws.Range("N134").Value = -31136.4

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 6420.857
$ws.Range("I136").Value = 3754.6
$ws.Range("J136").Value = 8844.727999999999
$ws.Range("K136").Value = 11263.8
$ws.Range("L136").Value = 26534.184
$ws.Range("M136").Value = -8713.799999999999
$ws.Range("N136").Value = -31634.184

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H26").Value = 105.14286
$ws.Range("J26").Value = 132
$ws.Range("L26").Value = 396
$ws.Range("N26").Value = -972

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 2812.724
$ws.Range("J34").Value = 6063.615
$ws.Range("L34").Value = 18190.845
$ws.Range("N34").Value = -18358.845

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H98").Value = 1635.1111
$ws.Range("I98").Value = 1594.5
$ws.Range("J98").Value = 1716.3334
$ws.Range("K98").Value = 4783.5
$ws.Range("L98").Value = 5149.0002
$ws.Range("M98").Value = -3285.5
$ws.Range("N98").Value = -8145.0002

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 24195
$ws.Range("I5").Value = 24195
$ws.Range("K5").Value = 24195
$ws.Range("M5").Value = -24083

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1011.2308
$ws.Range("I97").Value = 1021.625
$ws.Range("K97").Value = 1021.625
$ws.Range("M97").Value = -525.625

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 757.3077
$ws.Range("J107").Value = 608
$ws.Range("L107").Value = 608
$ws.Range("N107").Value = -4448

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 6000.1816
$ws.Range("I126").Value = 5111.8887
$ws.Range("J126").Value = 9997.5
$ws.Range("K126").Value = 15335.6661
$ws.Range("L126").Value = 29992.5
$ws.Range("M126").Value = -12865.6661
$ws.Range("N126").Value = -34932.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 24389.26
$ws.Range("I132").Value = 14977.125
$ws.Range("K132").Value = 44931.375
$ws.Range("M132").Value = -42401.375

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H12").Value = 0
$ws.Range("I12").Value = 0
$ws.Range("K12").Value = 0
$ws.Range("M12").Value = $null

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 3998.7222
$ws.Range("I93").Value = 5553
$ws.Range("J93").Value = 1556.2858
$ws.Range("K93").Value = 5553
$ws.Range("L93").Value = 1556.2858
$ws.Range("M93").Value = -4305
$ws.Range("N93").Value = -4052.2858

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 3251.7368
$ws.Range("I100").Value = 2892.2
$ws.Range("J100").Value = 4600
$ws.Range("K100").Value = 2892.2
$ws.Range("L100").Value = 4600
$ws.Range("M100").Value = -2351.2
$ws.Range("N100").Value = -5682

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 9366.429
$ws.Range("I132").Value = 9191.666999999999
$ws.Range("J132").Value = 9497.5
$ws.Range("K132").Value = 27575.001
$ws.Range("L132").Value = 28492.5
$ws.Range("M132").Value = -25045.001
$ws.Range("N132").Value = -33552.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 6879
$ws.Range("I136").Value = 6435.4443
$ws.Range("K136").Value = 19306.3329
$ws.Range("M136").Value = -16756.3329

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").Value = $null

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 2238.75
$ws.Range("I81").Value = 826.75
$ws.Range("J81").Value = 2591.75
$ws.Range("K81").Value = 1653.5
$ws.Range("L81").Value = 5183.5
$ws.Range("M81").Value = -592.5
$ws.Range("N81").Value = -7305.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H84").Value = 2238.75
$ws.Range("I84").Value = 826.75
$ws.Range("J84").Value = 2591.75
$ws.Range("K84").Value = 8267.5
$ws.Range("L84").Value = 25917.5
$ws.Range("M84").Value = -2963.5
$ws.Range("N84").Value = -36525.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 8548886
$ws.Range("I107").Value = 1405.2222
$ws.Range("J107").Value = 27780720
$ws.Range("K107").Value = 4215.6666
$ws.Range("L107").Value = 83342160
$ws.Range("M107").Value = -2295.6666
$ws.Range("N107").Value = -83346000

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H109").Value = 80890.10000000001
$ws.Range("I109").Value = 0
$ws.Range("K109").Value = 0
$ws.Range("M109").Value = $null

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H123").Value = 49999.547
$ws.Range("J123").Value = 49999.547
$ws.Range("L123").Value = 49999.547
$ws.Range("N123").Value = -59799.547

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 8002942.5
$ws.Range("I136").Value = 11113276
$ws.Range("K136").Value = 33339828
$ws.Range("M136").Value = -33337278

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H139").Value = 99999.664
$ws.Range("J139").Value = 99999.664
$ws.Range("L139").Value = 99999.664
$ws.Range("N139").Value = -110279.664
